# referrals.xlsx — "Put data together in gui class"
#
# The underlying text-file query (QueryTables/Connections) picked up a new
# leading DATE column and one more row of referral data; reproduce the
# resulting worksheet shape: insert a new column A, label it, append the new
# row, and grow the autofilter / defined names to match.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Make room for the new DATE column - shifts B:G (old A:F... originally
#    A:G) one column to the right, carrying values/styles/widths along.
$ws.Columns("A:A").Insert()

# New DATE column gets its own best-fit-ish width (short "yyyy:m:d:h" values).
$ws.Columns("A").ColumnWidth = 9

# 2) New column header
$ws.Range("A1").Value = "DATE"

# 3) The DEST_ZONE/NAME columns (B, F) pick up the same "stored as text"
#    formatting their neighbouring SENDER_NUMBER/NUMBER columns (C, G)
#    already had, across the whole table (rows 1-7). Doing this before
#    writing row 7 also keeps values like "+886910358944" / "0987017211"
#    as text instead of being coerced to numbers.
$ws.Range("B1:B7").NumberFormat = "@"
$ws.Range("C7").NumberFormat = "@"
$ws.Range("F1:F7").NumberFormat = "@"
$ws.Range("G7").NumberFormat = "@"

# 4) New referral row - ordered to match the shared-string insertion order
#    of the authored workbook.
$ws.Range("A7").Value = "2016:2:4:7"
$ws.Range("E7").Value = "TEST"
$ws.Range("F7").Value = "田凹凸"
$ws.Range("H7").Value = "一個測試"
$ws.Range("C7").Value = "+886910358944"
$ws.Range("G7").Value = "0987017211"
$ws.Range("D7").Value = "OFFICE_E"
# B7 stays blank (text-formatted empty cell, like B4).

# 5) Grow the AutoFilter to the new used range A1:H7. Re-applying over an
#    already-filtered region toggles it off, so cycle AutoFilterMode to make
#    sure it ends up ON with the new range.
$ws.AutoFilterMode = $false
[void]$ws.Range("A1:H7").AutoFilter()

# 6) Defined names follow the grown range too.
$wb.Names.Item("Sheet1!_FilterDatabase").RefersTo = "=Sheet1!`$A`$1:`$H`$7"
$wb.Names.Item("Sheet1!referrals").RefersTo = "=Sheet1!`$A`$1:`$H`$7"

# 7) Selection moved as part of the edit.
[void]$ws.Range("D13").Select()
